$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values for rows 2-10 (E,F,G,H,I,J,M,N,O,P,Q,R,S,T columns)
# A,B,C,D,K,L are unchanged by this update.

# Row 2
$ws.Range("E2").Value2 = 2
$ws.Range("F2").Value2 = 0.6666666666666666
$ws.Range("G2").Value2 = 0.2988413333333333
$ws.Range("H2").Value2 = 0.896524
$ws.Range("I2").Value2 = 0.3632971504731247
$ws.Range("J2").Value2 = 0.3632971504731246
$ws.Range("M2").Value2 = 2.211928
$ws.Range("N2").Value2 = 6.635783999999999
$ws.Range("O2").Value2 = 0.174938892641363
$ws.Range("P2").Value2 = 0.1749388926413629
$ws.Range("Q2").Value2 = 0.6610155127573333
$ws.Range("R2").Value2 = 5.949139614816
$ws.Range("S2").Value2 = 0.06355480120353105
$ws.Range("T2").Value2 = 0.06355480120353102

# Row 3
$ws.Range("E3").Value2 = 2
$ws.Range("F3").Value2 = 0.6666666666666666
$ws.Range("G3").Value2 = 0.2988413333333333
$ws.Range("H3").Value2 = 0.896524
$ws.Range("I3").Value2 = 0.3632971504731247
$ws.Range("J3").Value2 = 0.3632971504731246
$ws.Range("M3").Value2 = 3.864911333333334
$ws.Range("N3").Value2 = 11.594734
$ws.Range("O3").Value2 = 0.3056714815357404
$ws.Range("P3").Value2 = 0.3056714815357404
$ws.Range("Q3").Value2 = 1.154995256068445
$ws.Range("R3").Value2 = 10.394957304616
$ws.Range("S3").Value2 = 0.1110495782228328
$ws.Range("T3").Value2 = 0.1110495782228328

# Row 4
$ws.Range("E4").Value2 = 2
$ws.Range("F4").Value2 = 0.6666666666666666
$ws.Range("G4").Value2 = 0.2988413333333333
$ws.Range("H4").Value2 = 0.896524
$ws.Range("I4").Value2 = 0.3632971504731247
$ws.Range("J4").Value2 = 0.3632971504731246
$ws.Range("M4").Value2 = 6.567164333333333
$ws.Range("N4").Value2 = 19.701493
$ws.Range("O4").Value2 = 0.5193896258228967
$ws.Range("P4").Value2 = 0.5193896258228966
$ws.Range("Q4").Value2 = 1.962540145592444
$ws.Range("R4").Value2 = 17.662861310332
$ws.Range("S4").Value2 = 0.1886927710467608
$ws.Range("T4").Value2 = 0.1886927710467607

# Row 5
$ws.Range("E5").Value2 = 2
$ws.Range("F5").Value2 = 0.6666666666666666
$ws.Range("G5").Value2 = 0.3851916666666667
$ws.Range("H5").Value2 = 1.155575
$ws.Range("I5").Value2 = 0.4682720202225272
$ws.Range("J5").Value2 = 0.4682720202225272
$ws.Range("M5").Value2 = 2.211928
$ws.Range("N5").Value2 = 6.635783999999999
$ws.Range("O5").Value2 = 0.174938892641363
$ws.Range("P5").Value2 = 0.1749388926413629
$ws.Range("Q5").Value2 = 0.8520162328666666
$ws.Range("R5").Value2 = 7.668146095799999
$ws.Range("S5").Value2 = 0.08191898867266284
$ws.Range("T5").Value2 = 0.08191898867266283

# Row 6
$ws.Range("E6").Value2 = 2
$ws.Range("F6").Value2 = 0.6666666666666666
$ws.Range("G6").Value2 = 0.3851916666666667
$ws.Range("H6").Value2 = 1.155575
$ws.Range("I6").Value2 = 0.4682720202225272
$ws.Range("J6").Value2 = 0.4682720202225272
$ws.Range("M6").Value2 = 3.864911333333334
$ws.Range("N6").Value2 = 11.594734
$ws.Range("O6").Value2 = 0.3056714815357404
$ws.Range("P6").Value2 = 0.3056714815357404
$ws.Range("Q6").Value2 = 1.488731638005556
$ws.Range("R6").Value2 = 13.39858474205
$ws.Range("S6").Value2 = 0.1431374021831541
$ws.Range("T6").Value2 = 0.1431374021831541

# Row 7
$ws.Range("E7").Value2 = 2
$ws.Range("F7").Value2 = 0.6666666666666666
$ws.Range("G7").Value2 = 0.3851916666666667
$ws.Range("H7").Value2 = 1.155575
$ws.Range("I7").Value2 = 0.4682720202225272
$ws.Range("J7").Value2 = 0.4682720202225272
$ws.Range("M7").Value2 = 6.567164333333333
$ws.Range("N7").Value2 = 19.701493
$ws.Range("O7").Value2 = 0.5193896258228967
$ws.Range("P7").Value2 = 0.5193896258228966
$ws.Range("Q7").Value2 = 2.529616974830556
$ws.Range("R7").Value2 = 22.766552773475
$ws.Range("S7").Value2 = 0.2432156293667103
$ws.Range("T7").Value2 = 0.2432156293667103

# Row 8
$ws.Range("E8").Value2 = 2
$ws.Range("F8").Value2 = 0.6666666666666666
$ws.Range("G8").Value2 = 0.138548
$ws.Range("H8").Value2 = 0.415644
$ws.Range("I8").Value2 = 0.1684308293043481
$ws.Range("J8").Value2 = 0.1684308293043481
$ws.Range("M8").Value2 = 2.211928
$ws.Range("N8").Value2 = 6.635783999999999
$ws.Range("O8").Value2 = 0.174938892641363
$ws.Range("P8").Value2 = 0.1749388926413629
$ws.Range("Q8").Value2 = 0.306458200544
$ws.Range("R8").Value2 = 2.758123804896
$ws.Range("S8").Value2 = 0.02946510276516909
$ws.Range("T8").Value2 = 0.02946510276516909

# Row 9
$ws.Range("E9").Value2 = 2
$ws.Range("F9").Value2 = 0.6666666666666666
$ws.Range("G9").Value2 = 0.138548
$ws.Range("H9").Value2 = 0.415644
$ws.Range("I9").Value2 = 0.1684308293043481
$ws.Range("J9").Value2 = 0.1684308293043481
$ws.Range("M9").Value2 = 3.864911333333334
$ws.Range("N9").Value2 = 11.594734
$ws.Range("O9").Value2 = 0.3056714815357404
$ws.Range("P9").Value2 = 0.3056714815357404
$ws.Range("Q9").Value2 = 0.5354757354106667
$ws.Range("R9").Value2 = 4.819281618696
$ws.Range("S9").Value2 = 0.05148450112975349
$ws.Range("T9").Value2 = 0.05148450112975349

# Row 10
$ws.Range("E10").Value2 = 2
$ws.Range("F10").Value2 = 0.6666666666666666
$ws.Range("G10").Value2 = 0.138548
$ws.Range("H10").Value2 = 0.415644
$ws.Range("I10").Value2 = 0.1684308293043481
$ws.Range("J10").Value2 = 0.1684308293043481
$ws.Range("M10").Value2 = 6.567164333333333
$ws.Range("N10").Value2 = 19.701493
$ws.Range("O10").Value2 = 0.5193896258228967
$ws.Range("P10").Value2 = 0.5193896258228966
$ws.Range("Q10").Value2 = 0.9098674840546667
$ws.Range("R10").Value2 = 8.188807356491999
$ws.Range("S10").Value2 = 0.08748122540942556
$ws.Range("T10").Value2 = 0.08748122540942554

